$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume snapshot refresh (GitHub Actions bot).
# Column D (Price) cells are plain text in the source sheet even when
# they look numeric (e.g. "1.00", "304.51"); cells whose new value would
# otherwise be auto-parsed as a Number are pre-formatted as Text ("@")
# so .Value stores the exact original string instead of a coerced number.

$ws.Range("D2").Value = "45.902.76"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "2.582.68"
$ws.Range("E3").Value = "  +8.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.51"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.99"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +4.48%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.564"
$ws.Range("E9").Value = "  +10.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.90"
$ws.Range("E10").Value = "  +9.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0829"
$ws.Range("E11").Value = "  +5.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.95"
$ws.Range("E12").Value = "  +11.29%  "
$ws.Range("D13").Value = "2.972.58"
$ws.Range("E13").Value = "  +8.55%  "
$ws.Range("E14").Value = "  +0.89%  "
$ws.Range("D15").Value = "2.583.76"
$ws.Range("E15").Value = "  +9.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.889"
$ws.Range("E16").Value = "  +7.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.67"
$ws.Range("E17").Value = "  +6.63%  "
$ws.Range("D18").Value = "45.960.06"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  +1.67%  "
$ws.Range("D20").Value = "0.0₃0998"
$ws.Range("E20").Value = "  +4.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.59"
$ws.Range("E21").Value = "  +8.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.26"
$ws.Range("E22").Value = "  +5.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "251.90"
$ws.Range("E23").Value = "  +3.24%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.96"
$ws.Range("E24").Value = "  +5.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.19"
$ws.Range("E25").Value = "  +13.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.42"
$ws.Range("E26").Value = "  +30.81%  "
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("E28").Value = "  +5.59%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "39.06"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.25"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.05"
$ws.Range("E31").Value = "  +9.16%  "
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.93"
$ws.Range("E33").Value = "  +4.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.26"
$ws.Range("E34").Value = "  +17.45%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "152.05"
$ws.Range("E35").Value = "  +3.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0822"
$ws.Range("E36").Value = "  +6.12%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("E38").Value = "  +4.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.13"
$ws.Range("E39").Value = "  +5.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.60"
$ws.Range("E40").Value = "  +3.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.53"
$ws.Range("E41").Value = "  +9.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0318"
$ws.Range("E42").Value = "  +6.00%  "
$ws.Range("D43").Value = "2.044.34"
$ws.Range("E43").Value = "  +5.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.42"
$ws.Range("E44").Value = "  +35.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.998"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.61"
$ws.Range("E46").Value = "  -1.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.18"
$ws.Range("E47").Value = "  +7.88%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.76"
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "108.28"
$ws.Range("E49").Value = "  +9.03%  "
$ws.Range("D50").Value = "2.830.89"
$ws.Range("E50").Value = "  +8.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.198"
$ws.Range("E51").Value = "  +5.82%  "
